$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.371.54"
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").Value = "'1.565.33"
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = "'288.84"
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = "'0.3737"
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("D8").Value = "'49.22"
$ws.Range("D9").Value = "'0.3367"
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("D10").Value = "'0.07425"
$ws.Range("E10").Value = '  -2.54%  '
$ws.Range("D11").Value = "'1.117"
$ws.Range("E11").Value = '  -3.69%  '
$ws.Range("D12").Value = "'1.006"
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").Value = "'20.69"
$ws.Range("E13").Value = '  -3.18%  '
$ws.Range("D14").Value = "'5.868"
$ws.Range("E14").Value = '  -2.58%  '
$ws.Range("D15").Value = "'6.844"
$ws.Range("E15").Value = '  -0.76%  '
$ws.Range("D16").Value = "'1.562.27"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = "'0.00001107"
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").Value = "'89.09"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = "'0.06697"
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").Value = "'6.124"
$ws.Range("E21").Value = '  -1.32%  '
$ws.Range("D22").Value = "'16.19"
$ws.Range("E22").Value = '  -1.80%  '
$ws.Range("D23").Value = "'11.80"
$ws.Range("E23").Value = '  -1.23%  '
$ws.Range("D24").Value = "'22.369.83"
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").Value = "'2.366"
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("D26").Value = "'2.515"
$ws.Range("E26").Value = '  -10.55%  '
$ws.Range("D27").Value = "'19.87"
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("D28").Value = "'146.68"
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("D29").Value = "'4.990"
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").Value = "'124.71"
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("D31").Value = "'1.734.40"
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").Value = "'1.987"
$ws.Range("E32").Value = '  -1.00%  '
$ws.Range("E33").Value = '  -2.25%  '
$ws.Range("D34").Value = "'5.886"
$ws.Range("E34").Value = '  -4.29%  '
$ws.Range("D35").Value = "'9.700"
$ws.Range("E35").Value = '  -3.22%  '
$ws.Range("D36").Value = "'0.08411"
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("D37").Value = "'1.382"
$ws.Range("E37").Value = '  +3.13%  '
$ws.Range("D38").Value = "'0.02448"
$ws.Range("E38").Value = '  -3.48%  '
$ws.Range("D39").Value = "'0.2252"
$ws.Range("E39").Value = '  -2.36%  '
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").Value = "'5.318"
$ws.Range("E41").Value = '  -3.13%  '
$ws.Range("D42").Value = "'0.6161"
$ws.Range("E42").Value = '  -2.62%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = "'10.92"
$ws.Range("E43").Value = '  -6.26%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = "'1.005"
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("D46").Value = "'3.771"
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").Value = "'0.5742"
$ws.Range("E47").Value = '  -3.59%  '
$ws.Range("D48").Value = "'2.031"
$ws.Range("E48").Value = '  -2.55%  '
$ws.Range("D49").Value = "'124.65"
$ws.Range("E49").Value = '  +0.73%  '
$ws.Range("D50").Value = "'1.224"
$ws.Range("E50").Value = '  -2.76%  '
$ws.Range("D51").Value = "'0.07295"
$ws.Range("E51").Value = '  +0.59%  '
